$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new tariff row data
$ws.Range("A3").Value = 12345678901122
$ws.Range("B3").Value = 1860
$ws.Range("C3").Value = 9619

# Update the searchbox (active selection) to B3
$ws.Range("B3").Select()
